$wb = $excel.ActiveWorkbook

# --- News sheet: move selection to C3 (content already present) ---
$newsSheet = $wb.Worksheets.Item("News")
$newsSheet.Range("C3").Select() | Out-Null

# --- Add new sheets after "News": DeliveryTime, then Categories ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$deliverySheet = $wb.Worksheets.Add($null, $lastSheet)
$deliverySheet.Name = "DeliveryTime"

$catSheet = $wb.Worksheets.Add($null, $deliverySheet)
$catSheet.Name = "Categories"

# --- DeliveryTime sheet content ---
$deliverySheet.Range("A1").Value = 0.40625
$deliverySheet.Range("A1").NumberFormat = "hh:mm:ss"
$deliverySheet.Range("C5").Value = "  "
$deliverySheet.Range("C5").Select() | Out-Null

# --- Categories sheet content ---
$catSheet.Range("A1").Value = "Food"
$catSheet.Range("A2").Value = "Electronics"
$catSheet.Range("A3").Value = "Home Appliances"
$catSheet.Range("A4").Value = "Beauty"
$catSheet.Range("A4").Select() | Out-Null

# Categories ends up the active/selected sheet
$catSheet.Activate() | Out-Null
